# Auto-generated edit script applying the "Horarios actualizados Linea 141 - 614" update
# to the LP1912 / LP1912-215 / 6203-6173 schedule-scrape workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912": header + individual cell corrections, plus 7 newly scraped rows (350-356) ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 20:11:58"
$ws1.Cells.Item(3,1).Value = "Total filas: 351"
$ws1.Cells.Item(38,3).Value = "15_ABASTO"
$ws1.Cells.Item(39,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(49,1).Value = "08:52:40"
$ws1.Cells.Item(49,3).Value = "215B_EL PATO"
$ws1.Cells.Item(49,4).Value = 0
$ws1.Cells.Item(50,1).Value = "07:13:03"
$ws1.Cells.Item(50,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(50,4).Value = 99
$ws1.Cells.Item(137,1).Value = "10:36:50"
$ws1.Cells.Item(137,3).Value = "15_ABASTO"
$ws1.Cells.Item(137,4).Value = 118
$ws1.Cells.Item(138,1).Value = "11:46:32"
$ws1.Cells.Item(138,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(138,4).Value = 48
$ws1.Cells.Item(147,1).Value = "11:33:52"
$ws1.Cells.Item(147,3).Value = "215C_EL PATO"
$ws1.Cells.Item(147,4).Value = 90
$ws1.Cells.Item(148,1).Value = "11:13:15"
$ws1.Cells.Item(148,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(148,4).Value = 110
$ws1.Cells.Item(158,1).Value = "11:53:44"
$ws1.Cells.Item(158,3).Value = "215A_EL PATO"
$ws1.Cells.Item(158,4).Value = 99
$ws1.Cells.Item(159,1).Value = "12:11:21"
$ws1.Cells.Item(159,3).Value = "14_ABASTO"
$ws1.Cells.Item(159,4).Value = 81
$ws1.Cells.Item(169,1).Value = "12:33:02"
$ws1.Cells.Item(169,3).Value = "10_OLMOS"
$ws1.Cells.Item(169,4).Value = 89
$ws1.Cells.Item(170,1).Value = "12:46:07"
$ws1.Cells.Item(170,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(170,4).Value = 76
$ws1.Cells.Item(171,1).Value = "13:14:31"
$ws1.Cells.Item(171,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(171,4).Value = 48
$ws1.Cells.Item(178,1).Value = "12:53:26"
$ws1.Cells.Item(178,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(178,4).Value = 84
$ws1.Cells.Item(179,1).Value = "12:33:02"
$ws1.Cells.Item(179,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(179,4).Value = 104
$ws1.Cells.Item(235,1).Value = "16:44:58"
$ws1.Cells.Item(235,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(235,4).Value = 22
$ws1.Cells.Item(236,1).Value = "16:28:21"
$ws1.Cells.Item(236,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(236,4).Value = 38
$ws1.Cells.Item(247,1).Value = "17:35:41"
$ws1.Cells.Item(247,3).Value = "215B_EL PATO"
$ws1.Cells.Item(247,4).Value = 2
$ws1.Cells.Item(248,1).Value = "16:12:06"
$ws1.Cells.Item(248,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(248,4).Value = 85
$ws1.Cells.Item(258,1).Value = "16:51:51"
$ws1.Cells.Item(258,3).Value = "10_OLMOS"
$ws1.Cells.Item(258,4).Value = 62
$ws1.Cells.Item(259,1).Value = "16:37:37"
$ws1.Cells.Item(259,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(259,4).Value = 76
$ws1.Cells.Item(267,3).Value = "15_ABASTO"
$ws1.Cells.Item(268,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(269,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(270,3).Value = "15_ABASTO"
$ws1.Cells.Item(301,1).Value = "17:55:25"
$ws1.Cells.Item(301,3).Value = "17_ROMERO"
$ws1.Cells.Item(301,4).Value = 81
$ws1.Cells.Item(302,1).Value = "17:35:41"
$ws1.Cells.Item(302,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(302,4).Value = 101
$ws1.Cells.Item(327,1).Value = "20:11:58"
$ws1.Cells.Item(327,2).Value = "20:12"
$ws1.Cells.Item(327,3).Value = "10_OLMOS"
$ws1.Cells.Item(327,4).Value = 1
$ws1.Cells.Item(328,1).Value = "19:11:44"
$ws1.Cells.Item(328,2).Value = "20:13"
$ws1.Cells.Item(328,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(328,4).Value = 62
$ws1.Cells.Item(329,1).Value = "18:30:48"
$ws1.Cells.Item(329,2).Value = "20:21"
$ws1.Cells.Item(329,4).Value = 111
$ws1.Cells.Item(330,1).Value = "18:52:29"
$ws1.Cells.Item(330,2).Value = "20:22"
$ws1.Cells.Item(330,3).Value = "15_ABASTO"
$ws1.Cells.Item(330,4).Value = 90
$ws1.Cells.Item(331,1).Value = "18:44:45"
$ws1.Cells.Item(331,2).Value = "20:30"
$ws1.Cells.Item(331,3).Value = "10_OLMOS"
$ws1.Cells.Item(331,4).Value = 106
$ws1.Cells.Item(332,1).Value = "19:47:50"
$ws1.Cells.Item(332,2).Value = "20:33"
$ws1.Cells.Item(332,4).Value = 46
$ws1.Cells.Item(333,1).Value = "19:35:34"
$ws1.Cells.Item(333,2).Value = "20:34"
$ws1.Cells.Item(333,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(333,4).Value = 59
$ws1.Cells.Item(334,1).Value = "19:11:44"
$ws1.Cells.Item(334,2).Value = "20:41"
$ws1.Cells.Item(334,4).Value = 90
$ws1.Cells.Item(335,1).Value = "18:52:29"
$ws1.Cells.Item(335,2).Value = "20:42"
$ws1.Cells.Item(335,4).Value = 110
$ws1.Cells.Item(336,1).Value = "19:35:34"
$ws1.Cells.Item(336,2).Value = "20:43"
$ws1.Cells.Item(336,4).Value = 68
$ws1.Cells.Item(337,1).Value = "19:47:50"
$ws1.Cells.Item(337,2).Value = "20:45"
$ws1.Cells.Item(337,3).Value = "17_ROMERO"
$ws1.Cells.Item(337,4).Value = 58
$ws1.Cells.Item(338,1).Value = "18:52:29"
$ws1.Cells.Item(338,2).Value = "20:47"
$ws1.Cells.Item(338,3).Value = "215B_EL PATO"
$ws1.Cells.Item(338,4).Value = 115
$ws1.Cells.Item(339,1).Value = "20:11:58"
$ws1.Cells.Item(339,2).Value = "20:54"
$ws1.Cells.Item(339,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(339,4).Value = 43
$ws1.Cells.Item(340,1).Value = "19:35:34"
$ws1.Cells.Item(340,2).Value = "20:55"
$ws1.Cells.Item(340,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(340,4).Value = 80
$ws1.Cells.Item(341,1).Value = "19:54:57"
$ws1.Cells.Item(341,2).Value = "20:55"
$ws1.Cells.Item(341,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(341,4).Value = 61
$ws1.Cells.Item(342,1).Value = "19:11:44"
$ws1.Cells.Item(342,2).Value = "20:56"
$ws1.Cells.Item(342,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(342,4).Value = 105
$ws1.Cells.Item(343,1).Value = "19:11:44"
$ws1.Cells.Item(343,2).Value = "21:06"
$ws1.Cells.Item(343,3).Value = "10_OLMOS"
$ws1.Cells.Item(343,4).Value = 115
$ws1.Cells.Item(344,1).Value = "20:11:58"
$ws1.Cells.Item(344,2).Value = "21:06"
$ws1.Cells.Item(344,3).Value = "14_ABASTO"
$ws1.Cells.Item(344,4).Value = 55
$ws1.Cells.Item(345,2).Value = "21:09"
$ws1.Cells.Item(345,3).Value = "15_ABASTO"
$ws1.Cells.Item(345,4).Value = 82
$ws1.Cells.Item(346,1).Value = "19:35:34"
$ws1.Cells.Item(346,2).Value = "21:10"
$ws1.Cells.Item(346,3).Value = "15_ABASTO"
$ws1.Cells.Item(346,4).Value = 95
$ws1.Cells.Item(347,1).Value = "20:11:58"
$ws1.Cells.Item(347,2).Value = "21:27"
$ws1.Cells.Item(347,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(347,4).Value = 76
$ws1.Cells.Item(348,1).Value = "19:35:34"
$ws1.Cells.Item(348,2).Value = "21:28"
$ws1.Cells.Item(348,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(348,4).Value = 113
$ws1.Cells.Item(349,1).Value = "19:54:57"
$ws1.Cells.Item(349,2).Value = "21:33"
$ws1.Cells.Item(349,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(349,4).Value = 99
$ws1.Cells.Item(350,1).Value = "19:47:50"
$ws1.Cells.Item(350,2).Value = "21:33"
$ws1.Cells.Item(350,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(350,4).Value = 106
$ws1.Cells.Item(350,5).Value = "LP1912"
$ws1.Cells.Item(351,1).Value = "19:35:34"
$ws1.Cells.Item(351,2).Value = "21:34"
$ws1.Cells.Item(351,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(351,4).Value = 119
$ws1.Cells.Item(351,5).Value = "LP1912"
$ws1.Cells.Item(352,1).Value = "20:11:58"
$ws1.Cells.Item(352,2).Value = "21:37"
$ws1.Cells.Item(352,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(352,4).Value = 86
$ws1.Cells.Item(352,5).Value = "LP1912"
$ws1.Cells.Item(353,1).Value = "19:54:57"
$ws1.Cells.Item(353,2).Value = "21:44"
$ws1.Cells.Item(353,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(353,4).Value = 110
$ws1.Cells.Item(353,5).Value = "LP1912"
$ws1.Cells.Item(354,1).Value = "19:47:50"
$ws1.Cells.Item(354,2).Value = "21:45"
$ws1.Cells.Item(354,3).Value = "14X44_ABASTO"
$ws1.Cells.Item(354,4).Value = 118
$ws1.Cells.Item(354,5).Value = "LP1912"
$ws1.Cells.Item(355,1).Value = "20:11:58"
$ws1.Cells.Item(355,2).Value = "21:48"
$ws1.Cells.Item(355,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(355,4).Value = 97
$ws1.Cells.Item(355,5).Value = "LP1912"
$ws1.Cells.Item(356,1).Value = "20:11:58"
$ws1.Cells.Item(356,2).Value = "22:03"
$ws1.Cells.Item(356,3).Value = "15_ABASTO"
$ws1.Cells.Item(356,4).Value = 112
$ws1.Cells.Item(356,5).Value = "LP1912"

# --- Sheet "LP1912-215": refresh timestamp only ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 20:11:58"

# --- Sheet "6203-6173": refresh timestamp only ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 20:11:58"

